# Apply updated market-price / profit figures across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4233.778
$ws.Range("I64").Value = 4864.6665
$ws.Range("J64").Value = 3918.3333
$ws.Range("K64").Value = 4864.6665
$ws.Range("L64").Value = 3918.3333
$ws.Range("M64").Value = -4616.6665
$ws.Range("N64").Value = -4414.3333

$ws.Range("H67").Value = 4233.778
$ws.Range("I67").Value = 4864.6665
$ws.Range("J67").Value = 3918.3333
$ws.Range("K67").Value = 4864.6665
$ws.Range("L67").Value = 3918.3333
$ws.Range("M67").Value = -4006.6665
$ws.Range("N67").Value = -5634.3333

$ws.Range("H113").Value = 3650.6
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 3750.75
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 3750.75
$ws.Range("M113").Value = 4
$ws.Range("N113").Value = -10258.75

$ws.Range("H132").Value = 15882393
$ws.Range("I132").Value = 17553798
$ws.Range("J132").Value = 4050
$ws.Range("K132").Value = 52661394
$ws.Range("L132").Value = 12150
$ws.Range("M132").Value = -52658864
$ws.Range("N132").Value = -17210

$ws.Range("H138").Value = 1329.94
$ws.Range("I138").Value = 863.3077
$ws.Range("J138").Value = 1628.2787
$ws.Range("K138").Value = 2589.9231
$ws.Range("L138").Value = 4884.8361
$ws.Range("M138").Value = 2550.0769
$ws.Range("N138").Value = -15164.8361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 2788.1667
$ws.Range("I28").Value = 2788.1667
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2788.1667
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -2596.1667

$ws.Range("H32").Value = 4242.643
$ws.Range("I32").Value = 3850.6924
$ws.Range("J32").Value = 9338
$ws.Range("K32").Value = 3850.6924
$ws.Range("L32").Value = 9338
$ws.Range("M32").Value = -3563.6924

$ws.Range("H74").Value = 1797.8
$ws.Range("I74").Value = 1006
$ws.Range("J74").Value = 1995.75
$ws.Range("K74").Value = 1006
$ws.Range("L74").Value = 1995.75
$ws.Range("M74").Value = -132

$ws.Range("H77").Value = 1797.8
$ws.Range("I77").Value = 1006
$ws.Range("J77").Value = 1995.75
$ws.Range("K77").Value = 5030
$ws.Range("L77").Value = 9978.75
$ws.Range("M77").Value = -662

$ws.Range("H99").Value = 2788.1667
$ws.Range("I99").Value = 2788.1667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2788.1667
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 206.8332999999998

$ws.Range("H110").Value = 1269.7368
$ws.Range("I110").Value = 820.7778
$ws.Range("J110").Value = 1673.8
$ws.Range("K110").Value = 820.7778
$ws.Range("L110").Value = 1673.8
$ws.Range("M110").Value = 1224.2222
$ws.Range("N110").Value = -5763.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 16667555
$ws.Range("I94").Value = 20834144
$ws.Range("J94").Value = 1196.3334
$ws.Range("K94").Value = 20834144
$ws.Range("L94").Value = 1196.3334
$ws.Range("M94").Value = -20833693

$ws.Range("H134").Value = 12395.917
$ws.Range("I134").Value = 1350.4
$ws.Range("J134").Value = 20285.572
$ws.Range("K134").Value = 4051.2
$ws.Range("L134").Value = 60856.716
$ws.Range("M134").Value = -1516.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1703.4
$ws.Range("I31").Value = 1559.3334
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1559.3334
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1264.3334
$ws.Range("N31").Value = -3590

$ws.Range("H34").Value = 1703.4
$ws.Range("I34").Value = 1559.3334
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1559.3334
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1357.3334
$ws.Range("N34").Value = -3404

$ws.Range("H58").Value = 1252.2106
$ws.Range("I58").Value = 1053.0385
$ws.Range("J58").Value = 1683.75
$ws.Range("K58").Value = 1053.0385
$ws.Range("L58").Value = 1683.75
$ws.Range("M58").Value = -850.0385000000001

$ws.Range("H99").Value = 2116.138
$ws.Range("I99").Value = 1826.8182
$ws.Range("J99").Value = 3025.4285
$ws.Range("K99").Value = 1826.8182
$ws.Range("L99").Value = 3025.4285
$ws.Range("M99").Value = -328.8181999999999
$ws.Range("N99").Value = -6021.4285

$ws.Range("H122").Value = 1057.85
$ws.Range("I122").Value = 911.86664
$ws.Range("J122").Value = 1495.8
$ws.Range("K122").Value = 2735.59992
$ws.Range("L122").Value = 4487.4
$ws.Range("M122").Value = -285.5999199999997

$ws.Range("H126").Value = 2116.138
$ws.Range("I126").Value = 1826.8182
$ws.Range("J126").Value = 3025.4285
$ws.Range("K126").Value = 5480.4546
$ws.Range("L126").Value = 9076.2855
$ws.Range("M126").Value = -3010.4546
$ws.Range("N126").Value = -14016.2855

$ws.Range("H134").Value = 1851.1875
$ws.Range("I134").Value = 1866
$ws.Range("J134").Value = 1798.2858
$ws.Range("K134").Value = 5598
$ws.Range("L134").Value = 5394.857400000001
$ws.Range("M134").Value = -3063

$ws.Range("H136").Value = 1252.2106
$ws.Range("I136").Value = 1053.0385
$ws.Range("J136").Value = 1683.75
$ws.Range("K136").Value = 3159.1155
$ws.Range("L136").Value = 5051.25
$ws.Range("M136").Value = -609.1155000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1175.258
$ws.Range("I5").Value = 1181.1
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 3543.3
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -3431.3

$ws.Range("H32").Value = 1517.5
$ws.Range("I32").Value = 702
$ws.Range("J32").Value = 2333
$ws.Range("K32").Value = 2106
$ws.Range("L32").Value = 6999
$ws.Range("M32").Value = -1823
$ws.Range("N32").Value = -7565

$ws.Range("H92").Value = 678.3
$ws.Range("I92").Value = 711.8570999999999
$ws.Range("J92").Value = 600
$ws.Range("K92").Value = 2135.5713
$ws.Range("L92").Value = 1800
$ws.Range("M92").Value = -887.5712999999996
$ws.Range("N92").Value = -4296

$ws.Range("H132").Value = 748.3333
$ws.Range("I132").Value = 742.6667
$ws.Range("J132").Value = 750.2222
$ws.Range("K132").Value = 6684.0003
$ws.Range("L132").Value = 6751.999800000001
$ws.Range("M132").Value = -4154.0003
$ws.Range("N132").Value = -11811.9998

$ws.Range("H135").Value = 1175.258
$ws.Range("I135").Value = 1181.1
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 10629.9
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -8094.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 75002200
$ws.Range("I70").Value = 83335070
$ws.Range("J70").Value = 66669332
$ws.Range("K70").Value = 83335070
$ws.Range("L70").Value = 66669332
$ws.Range("M70").Value = -83334800
$ws.Range("N70").Value = -66669872

$ws.Range("H73").Value = 75002200
$ws.Range("I73").Value = 83335070
$ws.Range("J73").Value = 66669332
$ws.Range("K73").Value = 83335070
$ws.Range("L73").Value = 66669332
$ws.Range("M73").Value = -83334134
$ws.Range("N73").Value = -66671204

$ws.Range("H102").Value = 1108.9
$ws.Range("I102").Value = 1779.5
$ws.Range("J102").Value = 773.6
$ws.Range("K102").Value = 1779.5
$ws.Range("L102").Value = 773.6
$ws.Range("M102").Value = -157.5

$ws.Range("H126").Value = 2394.8462
$ws.Range("I126").Value = 2348.111
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 7044.333
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -4574.333

$ws.Range("H130").Value = 34620
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 34620
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 34620
$ws.Range("N130").Value = -44660

$ws.Range("H131").Value = 33950
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 33950
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 33950
$ws.Range("N131").Value = -44030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1368.3529
$ws.Range("I61").Value = 1287.8462
$ws.Range("J61").Value = 1630
$ws.Range("K61").Value = 1287.8462
$ws.Range("L61").Value = 1630
$ws.Range("M61").Value = -1085.8462
$ws.Range("N61").Value = -2034

$ws.Range("H100").Value = 1198
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 1330
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 1330
$ws.Range("M100").Value = -459
$ws.Range("N100").Value = -2412

$ws.Range("H113").Value = 1368.3529
$ws.Range("I113").Value = 1287.8462
$ws.Range("J113").Value = 1630
$ws.Range("K113").Value = 1287.8462
$ws.Range("L113").Value = 1630
$ws.Range("M113").Value = 882.1538
$ws.Range("N113").Value = -5970

$ws.Range("H132").Value = 23977.844
$ws.Range("I132").Value = 1412.5416
$ws.Range("J132").Value = 49766.76
$ws.Range("K132").Value = 4237.6248
$ws.Range("L132").Value = 149300.28
$ws.Range("M132").Value = -1707.6248
$ws.Range("N132").Value = -154360.28

$ws.Range("H136").Value = 14061.625
$ws.Range("I136").Value = 14061.625
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 42184.875
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -39634.875
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4697
$ws.Range("I132").Value = 4230.1665
$ws.Range("J132").Value = 5817.4
$ws.Range("K132").Value = 12690.4995
$ws.Range("L132").Value = 17452.2
$ws.Range("M132").Value = -10160.4995
$ws.Range("N132").Value = -22512.2

$ws.Range("H136").Value = 608.7143
$ws.Range("I136").Value = 513.1177
$ws.Range("J136").Value = 1015
$ws.Range("K136").Value = 1539.3531
$ws.Range("L136").Value = 3045
$ws.Range("M136").Value = 1010.6469
$ws.Range("N136").Value = -8145
